$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New translation-pair block (rows 15-16), mirroring the existing rows' layout:
# A = script path / section label, B = line numbers, C = English, D = Russian
# translation, E = "converted" (font-mapped) Russian string.
$ws.Range("A15").Value = 'SCRIPT/T01P01A/us0310.ssb'
$ws.Range("B15").Value = '21, 25, 29'
$ws.Range("C15").Value = ' The beach isn\''t that way,\n[hero].'
$ws.Range("D15").Value = ' [hero], пляж не в той\nстороне.'
$ws.Range("E15").Value = ' [hero], ðìÿç îå â óïê\nòóïñïîå.'

$ws.Range("A16").Value = 'У вас в группе Манафи'

# Match the row styling used by the rest of the data rows (vertical-top,
# wrap-text style index 4) and the author's row height for these entries.
$ws.Rows.Item(15).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 28.8

# Scroll the view down and move the selection the way the author left it.
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
